# Upgrade vars.xlsx "outputs" sheet to version 1.11.2.2
# - The outputs table (Table32103 / sheet3 "outputs") is reorganized: the
#   bastion/webapp/db/VPN-peer/ssh output groups are split into individual
#   rows (one server per row instead of zone1/zone2 pairs sharing a row),
#   two "ansible" helper rows become "ansible-bastion" / "ansible-vpn", and
#   the table grows from 23 to 27 rows (A1:D23 -> A1:D27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outputs")

$xlLeft = -4131

# ---- helper: write one data row (A/B/C text, D left blank) -------------
# NOTE: use ${row} (braced) instead of $row before a literal ":" --
# "A$row:D$row" mis-parses (the interpolation swallows past the colon),
# while "A${row}:D${row}" expands correctly.
function Set-OutRow($row, $a, $b, $c) {
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("A${row}:D${row}").HorizontalAlignment = $xlLeft
}

function Clear-OutRow($row) {
    $ws.Range("A$row").Value = ""
    $ws.Range("B$row").Value = ""
    $ws.Range("C$row").Value = ""
    $ws.Range("D$row").Value = ""
    $ws.Range("A${row}:D${row}").HorizontalAlignment = $xlLeft
}

# Row 1 (header) is untouched: *file / *name / *value / comments

Set-OutRow 2  "outputs.tf" "app_name" '"${var.dns-name}${var.domain}"'

Clear-OutRow 3

Set-OutRow 4  "outputs.tf" "master_db"       "ibm_is_instance.dbserver-zone1[0].primary_network_interface[0].primary_ipv4_address"
Set-OutRow 5  "outputs.tf" "slave_db"        "ibm_is_instance.dbserver-zone2[0].primary_network_interface[0].primary_ipv4_address"

Clear-OutRow 6

Set-OutRow 7  "outputs.tf" "webappserver1"   "ibm_is_instance.webappserver-zone1[0].primary_network_interface[0].primary_ipv4_address"
Set-OutRow 8  "outputs.tf" "webappserver2"   "ibm_is_instance.webappserver-zone2[0].primary_network_interface[0].primary_ipv4_address"

Clear-OutRow 9

Set-OutRow 10 "outputs.tf" "bastionserver1"  "ibm_is_floating_ip.bastionserver-zone1-fip.address"
Set-OutRow 11 "outputs.tf" "bastionserver2"  "ibm_is_floating_ip.bastionserver-zone2-fip.address"

Clear-OutRow 12

Set-OutRow 13 "outputs.tf" "ssh-bastionserver1" '"ssh root@${ibm_is_floating_ip.bastionserver-zone1-fip.address}"'
Set-OutRow 14 "outputs.tf" "ssh-webappserver1" '"ssh -o ProxyJump=root@${ibm_is_floating_ip.bastionserver-zone1-fip.address} root@${ibm_is_instance.webappserver-zone1[0].primary_network_interface[0].primary_ipv4_address}"'
Set-OutRow 15 "outputs.tf" "ssh-masterdb" '"ssh -o ProxyJump=root@${ibm_is_floating_ip.bastionserver-zone1-fip.address} root@${ibm_is_instance.dbserver-zone1[0].primary_network_interface[0].primary_ipv4_address}"'

Clear-OutRow 16

Set-OutRow 17 "outputs.tf" "ssh-bastionserver2" '"ssh root@${ibm_is_floating_ip.bastionserver-zone2-fip.address}"'
Set-OutRow 18 "outputs.tf" "ssh-webappserver2" '"ssh -o ProxyJump=root@${ibm_is_floating_ip.bastionserver-zone2-fip.address} root@${ibm_is_instance.webappserver-zone2[0].primary_network_interface[0].primary_ipv4_address}"'
Set-OutRow 19 "outputs.tf" "ssh-slavedb" '"ssh -o ProxyJump=root@${ibm_is_floating_ip.bastionserver-zone2-fip.address} root@${ibm_is_instance.dbserver-zone2[0].primary_network_interface[0].primary_ipv4_address}"'

Clear-OutRow 20

Set-OutRow 21 "outputs.tf" "VPN-peer1" "ibm_is_vpn_gateway.VPNGateway1.public_ip_address"
Set-OutRow 22 "outputs.tf" "VPN-peer2" "ibm_is_vpn_gateway.VPNGateway2.public_ip_address"

Clear-OutRow 23

Set-OutRow 24 "outputs.tf" "ansible-vpn" '"ansible-playbook -i inventory site.yaml"'
Set-OutRow 25 "outputs.tf" "ansible-bastion" '"ansible-playbook -i inventory --ssh-extra-args=''-J root@${ibm_is_floating_ip.bastionserver-zone1-fip.address}'' site.yaml"'

Clear-OutRow 26

# trailing "totals-ish" marker row, just a pair of single spaces like before
$ws.Range("A27").Value = ""
$ws.Range("B27").Value = " "
$ws.Range("C27").Value = " "
$ws.Range("D27").Value = ""
$ws.Range("A27:D27").HorizontalAlignment = $xlLeft

# ---- grow the ListObject / table to cover the new range ----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D27"))

# ---- selection matches the edited workbook (B23 was last touched cell) -
$ws.Activate()
$ws.Range("B23").Select()
